{"js": "// Replace the division expressions in the practice-sheet table.\n// Each entry is [rowIndex, colIndex, oldText, newText] in document order,\n// matching the table's row-major cell order (only the 5 \"content\" rows\n// hold text; the rest are blank spacer rows).\nconst replacements = [\n  [0, 0, \"23\u00f78=\", \"98\u00f76=\"],\n  [0, 1, \"28\u00f72=\", \"69\u00f77=\"],\n  [0, 2, \"20\u00f77=\", \"43\u00f72=\"],\n  [0, 3, \"57\u00f78=\", \"66\u00f78=\"],\n  [0, 4, \"97\u00f74=\", \"33\u00f77=\"],\n  [4, 0, \"17\u00f75=\", \"69\u00f79=\"],\n  [4, 1, \"29\u00f73=\", \"72\u00f79=\"],\n  [4, 2, \"53\u00f79=\", \"51\u00f74=\"],\n  [4, 3, \"76\u00f76=\", \"98\u00f74=\"],\n  [4, 4, \"76\u00f73=\", \"67\u00f73=\"],\n  [8, 0, \"30\u00f74=\", \"48\u00f75=\"],\n  [8, 1, \"73\u00f77=\", \"66\u00f77=\"],\n  [8, 2, \"15\u00f75=\", \"94\u00f73=\"],\n  [8, 3, \"27\u00f72=\", \"85\u00f75=\"],\n  [8, 4, \"47\u00f75=\", \"85\u00f73=\"],\n  [12, 0, \"53\u00f76=\", \"37\u00f75=\"],\n  [12, 1, \"43\u00f76=\", \"99\u00f77=\"],\n  [12, 2, \"47\u00f75=\", \"83\u00f73=\"],\n  [12, 3, \"12\u00f77=\", \"99\u00f76=\"],\n  [12, 4, \"59\u00f72=\", \"26\u00f75=\"],\n  [16, 0, \"27\u00f78=\", \"96\u00f77=\"],\n  [16, 1, \"53\u00f76=\", \"50\u00f76=\"],\n  [16, 2, \"73\u00f77=\", \"10\u00f74=\"],\n  [16, 3, \"57\u00f73=\", \"52\u00f76=\"],\n  [16, 4, \"33\u00f77=\", \"72\u00f72=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const [rowIndex, colIndex, oldText, newText] of replacements) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division expressions in the practice-sheet table.\n# Each entry is (row, column, oldText, newText) using 1-based table\n# coordinates; only rows 1, 5, 9, 13, 17 hold expressions (the rest are\n# blank spacer rows).\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$replacements = @(\n    ,@(1, 1, \"23\u00f78=\", \"98\u00f76=\")\n    ,@(1, 2, \"28\u00f72=\", \"69\u00f77=\")\n    ,@(1, 3, \"20\u00f77=\", \"43\u00f72=\")\n    ,@(1, 4, \"57\u00f78=\", \"66\u00f78=\")\n    ,@(1, 5, \"97\u00f74=\", \"33\u00f77=\")\n    ,@(5, 1, \"17\u00f75=\", \"69\u00f79=\")\n    ,@(5, 2, \"29\u00f73=\", \"72\u00f79=\")\n    ,@(5, 3, \"53\u00f79=\", \"51\u00f74=\")\n    ,@(5, 4, \"76\u00f76=\", \"98\u00f74=\")\n    ,@(5, 5, \"76\u00f73=\", \"67\u00f73=\")\n    ,@(9, 1, \"30\u00f74=\", \"48\u00f75=\")\n    ,@(9, 2, \"73\u00f77=\", \"66\u00f77=\")\n    ,@(9, 3, \"15\u00f75=\", \"94\u00f73=\")\n    ,@(9, 4, \"27\u00f72=\", \"85\u00f75=\")\n    ,@(9, 5, \"47\u00f75=\", \"85\u00f73=\")\n    ,@(13, 1, \"53\u00f76=\", \"37\u00f75=\")\n    ,@(13, 2, \"43\u00f76=\", \"99\u00f77=\")\n    ,@(13, 3, \"47\u00f75=\", \"83\u00f73=\")\n    ,@(13, 4, \"12\u00f77=\", \"99\u00f76=\")\n    ,@(13, 5, \"59\u00f72=\", \"26\u00f75=\")\n    ,@(17, 1, \"27\u00f78=\", \"96\u00f77=\")\n    ,@(17, 2, \"53\u00f76=\", \"50\u00f76=\")\n    ,@(17, 3, \"73\u00f77=\", \"10\u00f74=\")\n    ,@(17, 4, \"57\u00f73=\", \"52\u00f76=\")\n    ,@(17, 5, \"33\u00f77=\", \"72\u00f72=\")\n)\n\nforeach ($entry in $replacements) {\n    $rowIndex = $entry[0]\n    $colIndex = $entry[1]\n    $oldText = $entry[2]\n    $newText = $entry[3]\n\n    $cellRange = $table.Cell($rowIndex, $colIndex).Range\n    # Re-anchor a fresh document range over just this cell's character\n    # span; Find scoped to the cell's own Range object searches the whole\n    # document instead of just the cell, which breaks on duplicate text.\n    $scoped = $d.Range($cellRange.Start, $cellRange.End)\n\n    $scoped.Find.ClearFormatting()\n    $scoped.Find.Replacement.ClearFormatting()\n    $scoped.Find.Execute($oldText, $true, $true, $false, $false, $false, $true, 0, $false, $newText, 1) | Out-Null\n}\n"}
